$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("types")

# Insert a new row at 535 with the new reserved keyword "GLOBALID"
# (alphabetically between GET/GO in the Reserved Keywords list).
$ws.Rows.Item(535).Insert()
$ws.Range("A535").Value = "GLOBALID"

# The inserted row shifts the "Reserved" and "Special" named ranges
# down by 14 rows (they lived below the insertion point).
$n = $wb.Names.Item("Reserved")
$n.RefersTo = "=types!`$A`$178:`$A`$809"
$n2 = $wb.Names.Item("Special")
$n2.RefersTo = "=types!`$A`$167:`$A`$175"
